# "files for last week"
# Updates the Weekly sheet with a new task (DB Logic/Android Logic / DB Translator
# Logic) and reflects the corresponding totals / new task row on the Summary sheet.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly")
$wsSummary = $wb.Worksheets.Item("Summary")

# ----------------------------------------------------------------------
# Weekly sheet: log additional hours for the current week (column M) and
# add two brand-new task rows at the bottom of the sheet.
# ----------------------------------------------------------------------
$wsWeekly.Range("M13").Value = 4
$wsWeekly.Range("M15").Value = 8
$wsWeekly.Range("M16").Value = 4

$wsWeekly.Range("B17").Value = "DB Logic/Android Logic"
$wsWeekly.Range("M17").Value = 8

$wsWeekly.Range("B18").Value = "DB Translator Logic"
$wsWeekly.Range("M18").Value = 8

# ----------------------------------------------------------------------
# Summary sheet: update a couple of actuals totals, then insert a new row
# for the matching "Android/DB logic" task (pushing the rest of the
# Testing/Design/Integration sections down by one row; formulas adjust
# automatically).
# ----------------------------------------------------------------------
$wsSummary.Range("F9").Value = 27
$wsSummary.Range("F12").Value = 16

$wsSummary.Rows("18:18").Insert()
$wsSummary.Range("B18").Value = "Android/DB logic"
$wsSummary.Range("F18").Value = 8
$wsSummary.Range("C18").Clear()

# ----------------------------------------------------------------------
# View state: Weekly used to be the active/selected tab - make Summary the
# active tab instead, matching the last-viewed selections of each sheet.
# ----------------------------------------------------------------------
$wsWeekly.Range("M15").Select() | Out-Null

$wsSummary.Activate() | Out-Null
$wsSummary.Range("G16").Select() | Out-Null
